$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, avoiding Excel's automatic
# numeric/date coercion of strings like "244.23" or "1.000".
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '30.141.61'
Set-TextValue $ws.Range("E2") '  -0.65%  '
Set-TextValue $ws.Range("D3") '1.911.95'
Set-TextValue $ws.Range("E3") '  -1.16%  '
Set-TextValue $ws.Range("D4") '1.000'
Set-TextValue $ws.Range("E4") '  -0.08%  '
Set-TextValue $ws.Range("D5") '0.7401'
Set-TextValue $ws.Range("E5") '  -2.18%  '
Set-TextValue $ws.Range("D6") '244.23'
Set-TextValue $ws.Range("E6") '  -0.27%  '
Set-TextValue $ws.Range("D7") '1.000'
Set-TextValue $ws.Range("E7") '  -0.06%  '
Set-TextValue $ws.Range("D8") '0.3130'
Set-TextValue $ws.Range("E8") '  -1.42%  '
Set-TextValue $ws.Range("D9") '26.91'
Set-TextValue $ws.Range("E9") '  -2.36%  '
Set-TextValue $ws.Range("D10") '0.06983'
Set-TextValue $ws.Range("E10") '  -0.27%  '
Set-TextValue $ws.Range("D11") '0.7810'
Set-TextValue $ws.Range("E11") '  +0.42%  '
Set-TextValue $ws.Range("D12") '0.07994'
Set-TextValue $ws.Range("E12") '  -0.29%  '
Set-TextValue $ws.Range("D13") '1.901.87'
Set-TextValue $ws.Range("E13") '  -1.63%  '
Set-TextValue $ws.Range("D14") '5.288'
Set-TextValue $ws.Range("E14") '  -1.08%  '
Set-TextValue $ws.Range("D15") '92.27'
Set-TextValue $ws.Range("E15") '  -2.34%  '
Set-TextValue $ws.Range("D16") '14.41'
Set-TextValue $ws.Range("E16") '  +0.24%  '
Set-TextValue $ws.Range("D17") '30.150.54'
Set-TextValue $ws.Range("E17") '  -0.67%  '
Set-TextValue $ws.Range("D18") '5.926'
Set-TextValue $ws.Range("E18") '  +3.22%  '
Set-TextValue $ws.Range("D19") '241.94'
Set-TextValue $ws.Range("E19") '  -4.54%  '
Set-TextValue $ws.Range("D20") '0.000007848'
Set-TextValue $ws.Range("E20") '  -0.97%  '
Set-TextValue $ws.Range("D21") '1.000'
Set-TextValue $ws.Range("E21") '  +0.01%  '
Set-TextValue $ws.Range("D22") '2.158.36'
Set-TextValue $ws.Range("E22") '  -1.06%  '
Set-TextValue $ws.Range("D23") '1.000'
Set-TextValue $ws.Range("D24") '7.244'
Set-TextValue $ws.Range("E24") '  +8.76%  '
Set-TextValue $ws.Range("D25") '9.457'
Set-TextValue $ws.Range("E25") '  -0.15%  '
Set-TextValue $ws.Range("D26") '168.25'
Set-TextValue $ws.Range("E26") '  +1.32%  '
Set-TextValue $ws.Range("E27") '  +0.50%  '
Set-TextValue $ws.Range("D28") '0.1285'
Set-TextValue $ws.Range("E28") '  -3.81%  '
Set-TextValue $ws.Range("D29") '2.073'
Set-TextValue $ws.Range("E29") '  -5.70%  '
Set-TextValue $ws.Range("D30") '1.359'
Set-TextValue $ws.Range("E30") '  -0.57%  '
Set-TextValue $ws.Range("D31") '1.549'
Set-TextValue $ws.Range("E31") '  +2.05%  '
Set-TextValue $ws.Range("D32") '4.353'
Set-TextValue $ws.Range("E32") '  -1.01%  '
Set-TextValue $ws.Range("D33") '4.111'
Set-TextValue $ws.Range("E33") '  -0.22%  '
Set-TextValue $ws.Range("D34") '0.05186'
Set-TextValue $ws.Range("E34") '  +0.71%  '
Set-TextValue $ws.Range("D35") '1.300'
Set-TextValue $ws.Range("E35") '  +1.36%  '
Set-TextValue $ws.Range("D36") '0.7527'
Set-TextValue $ws.Range("E36") '  +0.21%  '
Set-TextValue $ws.Range("D37") '2.730'
Set-TextValue $ws.Range("E37") '  -1.40%  '
Set-TextValue $ws.Range("D38") '0.01946'
Set-TextValue $ws.Range("E38") '  -0.64%  '
Set-TextValue $ws.Range("D39") '2.794'
Set-TextValue $ws.Range("E39") '  -0.26%  '
Set-TextValue $ws.Range("D40") '6.380'
Set-TextValue $ws.Range("E40") '  -0.66%  '
Set-TextValue $ws.Range("D41") '0.4527'
Set-TextValue $ws.Range("E41") '  +1.50%  '
Set-TextValue $ws.Range("D42") '75.23'
Set-TextValue $ws.Range("E42") '  -2.76%  '
Set-TextValue $ws.Range("D43") '1.968'
Set-TextValue $ws.Range("E43") '  +0.12%  '
Set-TextValue $ws.Range("D44") '7.881'
Set-TextValue $ws.Range("E44") '  +5.27%  '
Set-TextValue $ws.Range("E45") '  +0.05%  '
Set-TextValue $ws.Range("D46") '0.8386'
Set-TextValue $ws.Range("E46") '  +0.68%  '
Set-TextValue $ws.Range("D47") '9.954'
Set-TextValue $ws.Range("E47") '  +1.88%  '
Set-TextValue $ws.Range("D48") '101.76'
Set-TextValue $ws.Range("E48") '  +1.14%  '
Set-TextValue $ws.Range("D49") '37.25'
Set-TextValue $ws.Range("E49") '  -0.73%  '
Set-TextValue $ws.Range("D50") '2.058.26'
Set-TextValue $ws.Range("E50") '  -0.96%  '

# Row 51: coin changed from Maker to Algorand
Set-TextValue $ws.Range("B51") 'Algorand'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D51") '0.1198'
Set-TextValue $ws.Range("E51") '  +1.58%  '

